# Daily attendance processing - 2025-10-05 10:16:28
# Applies the attendance-session refresh to "Session Analysis Results":
#  - Updated session statistics (recorded/missing/pending counts, coverage %)
#  - Re-synced "Recorded By" grader lists (order reshuffled by the sync job)
#  - Newly recorded / re-recorded sessions (students counts, status, grader)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value (e.g. "30.1%") without Excel
# auto-converting it to a percentage number - format as Text first,
# then restore the original cell style (alignment/fill) via a format-only paste.
function Set-TextValue($cellRef, $text, $formatSourceRef) {
  $c = $ws.Range($cellRef)
  $c.NumberFormat = "@"
  $c.Value = $text
  $ws.Range($formatSourceRef).Copy()
  $c.PasteSpecial(-4122)  # xlPasteFormats
  $excel.Application.CutCopyMode = $false
}

# Helper: re-stripe a row to the fill used for a given status (Pending/
# Not Recorded/Recorded) by copying formats from a row that already has it.
function Set-RowStatusFormat($sourceRange, $targetRange) {
  $ws.Range($sourceRange).Copy()
  $ws.Range($targetRange).PasteSpecial(-4122)  # xlPasteFormats
  $excel.Application.CutCopyMode = $false
}

# Row 6
$ws.Range("L6").Value = 46

# Row 7
$ws.Range("G7").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("L7").Value = 10

# Row 8
$ws.Range("H8").Value = "159/217"
$ws.Range("L8").Value = 97

# Row 9
Set-TextValue "L9" "30.1%" "L4"

# Row 10
Set-TextValue "L10" "49.0%" "L4"

# Row 14
$ws.Range("H14").Value = "200/217"

# Row 15
Set-TextValue "S15" "70.1%" "M15"

# Row 16
Set-TextValue "S16" "59.3%" "M15"

# Row 17
$ws.Range("G17").Value = "nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"

# Row 18
$ws.Range("P18").Value = 1
$ws.Range("Q18").Value = 9

# Row 20
$ws.Range("O20").Value = 5
$ws.Range("P20").Value = 0
Set-TextValue "R20" "29.4%" "M15"
Set-TextValue "S20" "27.9%" "M15"

# Row 23
$ws.Range("P23").Value = 4
$ws.Range("Q23").Value = 10

# Row 24
$ws.Range("G24").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# Row 25
$ws.Range("H25").Value = "126/216"

# Row 34
$ws.Range("G34").Value = "nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"

# Row 35
$ws.Range("G35").Value = "Salma.hassan@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"

# Row 45
$ws.Range("G45").Value = "Rania.a.youssef@med.asu.edu.eg, System, mohamed.saleem@med.asu.edu.eg, backup@backdoor.com"

# Row 51
$ws.Range("G51").Value = "Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"

# Row 52
$ws.Range("G52").Value = "yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"

# Row 54
Set-RowStatusFormat "A3:I3" "A54:I54"
$ws.Range("I54").Value = "Not Recorded"

# Row 62
$ws.Range("G62").Value = "Rania.a.youssef@med.asu.edu.eg, System, mohamed.saleem@med.asu.edu.eg, backup@backdoor.com"

# Row 68
$ws.Range("G68").Value = "Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"

# Row 69
$ws.Range("G69").Value = "yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"

# Row 72
$ws.Range("G72").Value = "mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

# Row 76
$ws.Range("G76").Value = "mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# Row 83
$ws.Range("G83").Value = "Youstina.ibrahim@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg"

# Row 85
$ws.Range("G85").Value = "Aya_hamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"

# Row 98
$ws.Range("G98").Value = "user@user.com, afaf.abdallah@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"

# Row 99
$ws.Range("G99").Value = "user@user.com, Walaa.h.ghanima@med.asu.edu.eg"

# Row 100
Set-RowStatusFormat "A7:I7" "A100:I100"
$ws.Range("G100").Value = "nourhanmohamed@med.asu.edu.eg"
$ws.Range("H100").Value = "50/224"
$ws.Range("I100").Value = "Recorded"

# Row 102
$ws.Range("G102").Value = "Aya_hamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"

# Row 109
$ws.Range("G109").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# Row 116
$ws.Range("G116").Value = "nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg"

# Row 119
$ws.Range("G119").Value = "yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg"

# Row 126
$ws.Range("G126").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# Row 133
$ws.Range("G133").Value = "nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg"

# Row 136
$ws.Range("G136").Value = "yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg"

# Row 143
$ws.Range("G143").Value = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# Row 147
Set-RowStatusFormat "A3:I3" "A147:I147"
$ws.Range("I147").Value = "Not Recorded"

# Row 149
$ws.Range("G149").Value = "user@user.com, Walaa.h.ghanima@med.asu.edu.eg"

# Row 150
$ws.Range("G150").Value = "Youstina.ibrahim@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
